# Common: Updated default atomizers
#
# Adds a "dual" column (C) to the "atomizers" sheet and appends two new
# rows for the "Vapor Giant" vendor ("Extreme 2" and "Extreme").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "atomizers" sheet

function Set-Text($cell, [string]$text) {
    # Plain text values - safe to assign directly.
    $cell.Value = $text
}

function Set-TrueFalseText($cell, [string]$text) {
    # Excel auto-coerces bare "true"/"false" strings assigned via .Value
    # into real booleans. The source workbook stores them as literal text
    # (shared string "true"/"false"), so build the text through a formula
    # on a scratch cell and paste back as a value, which keeps it as text.
    $stage = $ws.Cells.Item(1000, 26)
    $stage.Formula = '="' + $text + '"'
    $stage.Copy()
    $cell.PasteSpecial(-4163) # xlPasteValues
    $stage.ClearContents()
    $excel.CutCopyMode = $false
}

# New header cell for column C.
Set-Text $ws.Cells.Item(1, 3) "dual"

# Existing rows (2-9) get "false" in the new column C.
for ($r = 2; $r -le 9; $r++) {
    Set-TrueFalseText $ws.Cells.Item($r, 3) "false"
}

# New row 10: Vapor Giant / Extreme 2 / true
Set-Text $ws.Cells.Item(10, 1) "Vapor Giant"
Set-Text $ws.Cells.Item(10, 2) "Extreme 2"
Set-TrueFalseText $ws.Cells.Item(10, 3) "true"

# New row 11: Vapor Giant / Extreme / true
Set-Text $ws.Cells.Item(11, 1) "Vapor Giant"
Set-Text $ws.Cells.Item(11, 2) "Extreme"
Set-TrueFalseText $ws.Cells.Item(11, 3) "true"

# Match the selection shown in the diff.
$ws.Activate()
$ws.Range("C10").Select() | Out-Null
